$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Progress report: CB point 2 ("2. di Bank Payment warna Jumlah Bayar ...")
# is now Done instead of in Progress.
$ws.Range("B18").Value = "Done"

# Reflect where the user left the cursor/scroll in the sheet before saving.
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
